# Generate Report for Handback
# Update the handoff/handback timestamps on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-25 10:54:40"
$wsZhCn.Range("E3").Value = "2016-03-25 10:54:40"
$wsZhCn.Range("H2").Value = "2016-03-25 10:55:29"
$wsZhCn.Range("H3").Value = "2016-03-25 10:55:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-25 10:54:50"
$wsDeDe.Range("E3").Value = "2016-03-25 10:54:50"
$wsDeDe.Range("H2").Value = "2016-03-25 10:55:44"
$wsDeDe.Range("H3").Value = "2016-03-25 10:55:44"
